# Reorganized unit label in frequency UI.
#
# The "Translation" sheet (Worksheets index 2) lists UI text rows with
# columns: B=Text ID, C=Typography Name, D=Alignment, E=Direction, F=GB text.
#
# Rows 206-213 were the "displayLabel" rows for the uHz frequency unit
# (Left aligned, showing the literal unit text "uHz"). Rows 214-245 were the
# equivalent rows for mHz / kHz / MHz / GHz (8 rows each).
#
# The edit reorganizes this: the uHz group (206-213) and the mHz group
# (214-221, but re-keyed onto the TEXT IDs formerly used by the kHz group)
# are kept, now Center-aligned and showing the "<value>" placeholder instead
# of the hard-coded unit text. The remaining groups (what used to be rows
# 222-245) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Rows 206-213: uHz group -> Center aligned, "<value>" placeholder ---
# Text ID (col B) and Typography Name (col C) stay the same; only the
# alignment (D) and displayed text (F) change.
for ($r = 206; $r -le 213; $r++) {
    $ws.Cells.Item($r, 4).Value = "Center"
    $ws.Cells.Item($r, 6).Value = "<value>"
}

# --- Rows 214-221: re-keyed onto the Text IDs formerly used by the kHz
# group (rows 222-229), now also Center aligned with the "<value>" text ---
$newTextIds = @(
    "SingleUseId265",
    "SingleUseId266",
    "SingleUseId267",
    "SingleUseId268",
    "SingleUseId269",
    "SingleUseId270",
    "SingleUseId271",
    "SingleUseId272"
)
for ($i = 0; $i -lt 8; $i++) {
    $r = 214 + $i
    $ws.Cells.Item($r, 2).Value = $newTextIds[$i]
    $ws.Cells.Item($r, 4).Value = "Center"
    $ws.Cells.Item($r, 6).Value = "<value>"
}

# --- Rows 222-245: old kHz / MHz / GHz groups are removed ---
$ws.Range("B222:F245").ClearContents()
